$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.487.46'
$ws.Range('E2').Value = '  +0.86%  '
$ws.Range('D3').Value = '1.794.08'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '226.61'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('E6').Value = '  +1.61%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '32.50'
$ws.Range('E8').Value = '  +1.80%  '
$ws.Range('D9').Value = '0.296'
$ws.Range('E9').Value = '  +1.05%  '
$ws.Range('E10').Value = '  +0.55%  '
$ws.Range('D11').Value = '0.0950'
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('D12').Value = '2.053.71'
$ws.Range('E12').Value = '  +0.21%  '
$ws.Range('D13').Value = '11.03'
$ws.Range('E13').Value = '  -1.05%  '
$ws.Range('D14').Value = '1.787.17'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').Value = '0.637'
$ws.Range('E15').Value = '  +2.31%  '
$ws.Range('D16').Value = '34.455.10'
$ws.Range('E16').Value = '  +0.95%  '
$ws.Range('E17').Value = '  +2.08%  '
$ws.Range('D18').Value = '68.80'
$ws.Range('E18').Value = '  +0.93%  '
$ws.Range('D19').Value = '247.07'
$ws.Range('E19').Value = '  +0.59%  '
$ws.Range('D20').Value = '0.0₃0799'
$ws.Range('E20').Value = '  +2.55%  '
$ws.Range('D21').Value = '11.22'
$ws.Range('E21').Value = '  +3.43%  '
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D24').Value = '2.08'
$ws.Range('E24').Value = '  +1.00%  '
$ws.Range('D25').Value = '164.19'
$ws.Range('E25').Value = '  +1.94%  '
$ws.Range('D26').Value = '7.25'
$ws.Range('E26').Value = '  +1.05%  '
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('E28').Value = '  +2.25%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '3.80'
$ws.Range('E30').Value = '  +3.60%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.23'
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('E32').Value = '  +0.35%  '
$ws.Range('E33').Value = '  +6.63%  '
$ws.Range('D34').Value = '1.82'
$ws.Range('E34').Value = '  +0.71%  '
$ws.Range('D35').Value = '1.438.82'
$ws.Range('E35').Value = '  -0.99%  '
$ws.Range('D36').Value = '2.58'
$ws.Range('E36').Value = '  +6.47%  '
$ws.Range('E37').Value = '  +2.83%  '
$ws.Range('D38').Value = '1.07'
$ws.Range('E38').Value = '  +2.83%  '
$ws.Range('D39').Value = '0.0191'
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('D40').Value = '84.60'
$ws.Range('E40').Value = '  +5.19%  '
$ws.Range('E41').Value = '  +1.42%  '
$ws.Range('E42').Value = '  +1.49%  '
$ws.Range('E43').Value = '  +1.99%  '
$ws.Range('D44').Value = '13.50'
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('E45').Value = '  +3.60%  '
$ws.Range('D46').Value = '6.13'
$ws.Range('E46').Value = '  +1.17%  '
$ws.Range('E47').Value = '  +0.27%  '
$ws.Range('D48').Value = '1.950.68'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').Value = '105.77'
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0132'
$ws.Range('E50').Value = '  -3.10%  '
$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  +0.02%  '
